$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 45.6773314456603
$ws.Range("C2").Value = 0.40873237147177988
$ws.Range("D2").Value = 32.042792113907446
$ws.Range("E2").Value = 0.28845639571382869
$ws.Range("F2").Value = 0.53708136787066885
$ws.Range("G2").Value = 0.50649466706193247
$ws.Range("H2").Value = 0.71154360428617136
$ws.Range("I2").Value = 0.84719231370516268

# Row 3
$ws.Range("B3").Value = 47.126364750062592
$ws.Range("C3").Value = 0.42169868977683511
$ws.Range("D3").Value = 33.086163228014279
$ws.Range("E3").Value = 0.30704823541900106
$ws.Range("F3").Value = 0.55411933319367324
$ws.Range("G3").Value = 0.52298704710119293
$ws.Range("H3").Value = 0.69295176458099894
$ws.Range("I3").Value = 0.83288058624949202

# Row 4
$ws.Range("B4").Value = 48.388797894883567
$ws.Range("C4").Value = 0.43299526242624886
$ws.Range("D4").Value = 35.178524982019304
$ws.Range("E4").Value = 0.32371915026352605
$ws.Range("F4").Value = 0.56896322399916677
$ws.Range("G4").Value = 0.55606063401585903
$ws.Range("H4").Value = 0.67628084973647395
$ws.Range("I4").Value = 0.83429447792900602
